$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# Insert a new first column, shifting DPMMETRICDATATYPE..SEQUENCE_ID right by one.
$ws.Columns.Item(1).EntireColumn.Insert()

# The new column A becomes "MEMBER"; its values are the old SEQUENCE_ID/ORDER
# sequence numbers, which (after the shift) now live in column M.
$ws.Range("A1").Value = "MEMBER"
$ws.Range("M2:M31").Copy()
$ws.Range("A2").PasteSpecial(-4163)

# Drop the now-duplicated old SEQUENCE_ID column (shifted to N), shifting the
# remaining columns left so the sheet stays at 13 columns (A1:M31).
$ws.Range("N1:N31").Delete(-4159)

# Match the new column widths: D keeps the old C width, a new default-ish E
# column is widened to 30, and F keeps the old E width.
$ws.Columns.Item(5).ColumnWidth = 29.1666666666667

# Update the saved selection/view state.
$ws.Range("D8").Select()
